$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 31 (old rows 31.. shift down by one)
$ws.Rows("31").Insert()

# Populate the newly inserted row 31 with the new tokenization case
$ws.Range("A31").Value2 = 'echo "$"'
$ws.Range("B31").Value2 = 'echo $'
$ws.Range("C31").Value2 = 2

# Append two new rows (385, 386) at the end of the data for the new cases
$ws.Range("A385").Value2 = 'echo $"$SHLVL"'
$ws.Range("B385").Value2 = 'echo 2'
$ws.Range("C385").Value2 = 5

$ws.Range("A386").Value2 = 'echo $''$SHLVL'''
$ws.Range("B386").Value2 = 'echo $SHLVL'
$ws.Range("C386").Value2 = 5

# Match the author's final selection/view position
$ws.Range("C32").Select()
